# promo overlay msrp tests added
#
# Lexus MSRP staging sheet update:
#  - several NX trims move from MY2020 to MY2021 (rows 35-42, 43, 44)
#  - LC Convertible/LC Convertible trims move from MY2020 to MY2021 (rows 70, 71)
#  - the old "9260 (SE) / LC Inspiration Series" row is replaced with the new
#    "9262 / LC Convertible" pricing (row 72)
#  - a brand-new "9262SE / LC Convertible Inspiration Series" trim is appended
#    as row 80
#  - selection cursor moves to the new entry row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Model year bump: 2020 -> 2021 -------------------------------------
$ws.Range("C35:C42").Value = 2021
$ws.Range("C43").Value = 2021
$ws.Range("C44").Value = 2021
$ws.Range("C70").Value = 2021
$ws.Range("C71").Value = 2021

# --- Row 72: "9260 (SE)" / "LC Inspiration Series" -> "9262" / "LC Convertible"
$ws.Range("A72").Value = 9262
$ws.Range("B72").Value = "LC Convertible"
$ws.Range("C72").Value = 2021
$ws.Range("D72").Value = 101000

# --- New row 80: "9262SE" / "LC Convertible Inspiration Series" --------
$ws.Range("A80").Value = "9262SE"
$ws.Range("B80").Value = "LC Convertible Inspiration Series"
$ws.Range("C80").Value = 2021
$ws.Range("D80").Value = 119800
$ws.Range("D80").NumberFormat = $ws.Range("D79").NumberFormat
$ws.Range("E80").Value = 1025
$ws.Range("E80").NumberFormat = $ws.Range("E79").NumberFormat

# --- Move the selection cursor to just past the new last row -----------
[void]$ws.Range("B81").Select()
